$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range('D2')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '40.029.77'
$cell.Style = $origStyle
$ws.Range('E2').Value = '  -2.88%  '
$cell = $ws.Range('D3')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '2.339.03'
$cell.Style = $origStyle
$ws.Range('E3').Value = '  -3.96%  '
$ws.Range('E4').Value = '  -0.02%  '
$cell = $ws.Range('D5')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '310.05'
$cell.Style = $origStyle
$ws.Range('E5').Value = '  -2.04%  '
$cell = $ws.Range('D6')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '85.42'
$cell.Style = $origStyle
$ws.Range('E6').Value = '  -3.89%  '
$ws.Range('E7').Value = '  -2.27%  '
$ws.Range('E8').Value = '  -0.02%  '
$cell = $ws.Range('D9')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '0.484'
$cell.Style = $origStyle
$ws.Range('E9').Value = '  -2.55%  '
$cell = $ws.Range('D10')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '0.0811'
$cell.Style = $origStyle
$ws.Range('E10').Value = '  -2.63%  '
$cell = $ws.Range('D11')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '29.99'
$cell.Style = $origStyle
$ws.Range('E11').Value = '  -6.88%  '
$ws.Range('E12').Value = '  +1.02%  '
$ws.Range('E13').Value = '  -4.43%  '
$ws.Range('E14').Value = '  -4.37%  '
$cell = $ws.Range('D15')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '14.77'
$cell.Style = $origStyle
$ws.Range('E15').Value = '  -5.93%  '
$cell = $ws.Range('D16')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '2.364.44'
$cell.Style = $origStyle
$ws.Range('E16').Value = '  -2.61%  '
$ws.Range('E17').Value = '  -2.29%  '
$cell = $ws.Range('D18')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '40.006.01'
$cell.Style = $origStyle
$ws.Range('E18').Value = '  -2.78%  '
$cell = $ws.Range('D19')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '0.0₃0904'
$cell.Style = $origStyle
$ws.Range('E19').Value = '  -2.16%  '
$ws.Range('E20').Value = '  -1.90%  '
$cell = $ws.Range('D21')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '67.93'
$cell.Style = $origStyle
$ws.Range('E21').Value = '  -5.61%  '
$ws.Range('E22').Value = '  -3.54%  '
$cell = $ws.Range('D23')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '235.08'
$cell.Style = $origStyle
$ws.Range('E23').Value = '  -0.22%  '
$ws.Range('E24').Value = '  -5.18%  '
$ws.Range('E25').Value = '  +0.21%  '
$cell = $ws.Range('D26')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '1.82'
$cell.Style = $origStyle
$ws.Range('E26').Value = '  -3.14%  '
$cell = $ws.Range('D27')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '23.33'
$cell.Style = $origStyle
$ws.Range('E27').Value = '  -2.79%  '
$cell = $ws.Range('D28')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '2.12'
$cell.Style = $origStyle
$ws.Range('E28').Value = '  -4.16%  '
$ws.Range('E29').Value = '  -2.76%  '
$cell = $ws.Range('D30')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '34.96'
$cell.Style = $origStyle
$ws.Range('E30').Value = '  +0.61%  '
$cell = $ws.Range('D31')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '153.11'
$cell.Style = $origStyle
$ws.Range('E31').Value = '  -2.50%  '
$ws.Range('E32').Value = '  -0.10%  '
$ws.Range('E33').Value = '  -2.92%  '
$ws.Range('E34').Value = '  -3.04%  '
$cell = $ws.Range('D35')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '0.0719'
$cell.Style = $origStyle
$ws.Range('E35').Value = '  -3.40%  '
$ws.Range('E36').Value = '  -0.73%  '
$cell = $ws.Range('D37')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '2.80'
$cell.Style = $origStyle
$ws.Range('E37').Value = '  -4.36%  '
$cell = $ws.Range('D38')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '0.0990'
$cell.Style = $origStyle
$ws.Range('E38').Value = '  -0.98%  '
$cell = $ws.Range('D39')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '15.58'
$cell.Style = $origStyle
$ws.Range('E39').Value = '  -6.30%  '
$ws.Range('E40').Value = '  -3.05%  '
$cell = $ws.Range('D41')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '3.87'
$cell.Style = $origStyle
$ws.Range('E41').Value = '  +0.50%  '
$cell = $ws.Range('D42')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '1.957.42'
$cell.Style = $origStyle
$ws.Range('E42').Value = '  -1.33%  '
$ws.Range('E43').Value = '  -4.37%  '
$ws.Range('E44').Value = '  -4.51%  '
$cell = $ws.Range('D45')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '17.59'
$cell.Style = $origStyle
$ws.Range('E45').Value = '  -4.04%  '
$cell = $ws.Range('D46')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '9.39'
$cell.Style = $origStyle
$ws.Range('E46').Value = '  -0.99%  '
$ws.Range('E47').Value = '  -5.82%  '
$cell = $ws.Range('D48')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '2.561.13'
$cell.Style = $origStyle
$ws.Range('E48').Value = '  -3.87%  '
$cell = $ws.Range('D49')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '92.73'
$cell.Style = $origStyle
$ws.Range('E49').Value = '  -2.70%  '
$cell = $ws.Range('D50')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '70.69'
$cell.Style = $origStyle
$ws.Range('E50').Value = '  -3.44%  '
$cell = $ws.Range('D51')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '50.94'
$cell.Style = $origStyle
$ws.Range('E51').Value = '  -1.45%  '
